# Weekly update: refresh Hortaliza / Femacal de La Calera - Esparragos price data
# Each data row (2-28) is updated in place with the new date / quality / volume / price
# figures for the latest week, per the "Fruta / hortaliza, semanal" refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = 44174
$ws.Cells.Item(2, 9).Value = 'Primera'
$ws.Cells.Item(2, 10).Value = 2800
$ws.Cells.Item(2, 11).Value = 1200
$ws.Cells.Item(2, 12).Value = 1250
$ws.Cells.Item(2, 13).Value = 1221
$ws.Cells.Item(2, 15).Value = 'Provincia de Quillota'
$ws.Cells.Item(2, 16).Value = 1221

# Row 3
$ws.Cells.Item(3, 4).Value = 44174
$ws.Cells.Item(3, 9).Value = 'Segunda'
$ws.Cells.Item(3, 10).Value = 1300
$ws.Cells.Item(3, 11).Value = 1000
$ws.Cells.Item(3, 12).Value = 1000
$ws.Cells.Item(3, 13).Value = 1000
$ws.Cells.Item(3, 15).Value = 'Provincia de Quillota'
$ws.Cells.Item(3, 16).Value = 1000

# Row 4
$ws.Cells.Item(4, 4).Value = 44165
$ws.Cells.Item(4, 9).Value = 'Primera'
$ws.Cells.Item(4, 10).Value = 720
$ws.Cells.Item(4, 11).Value = 1200
$ws.Cells.Item(4, 12).Value = 1200
$ws.Cells.Item(4, 13).Value = 1200
$ws.Cells.Item(4, 15).Value = 'Provincia de Quillota'
$ws.Cells.Item(4, 16).Value = 1200

# Row 5
$ws.Cells.Item(5, 4).Value = 44165
$ws.Cells.Item(5, 9).Value = 'Segunda'
$ws.Cells.Item(5, 10).Value = 750
$ws.Cells.Item(5, 11).Value = 1000
$ws.Cells.Item(5, 12).Value = 1000
$ws.Cells.Item(5, 13).Value = 1000
$ws.Cells.Item(5, 15).Value = 'Provincia de Quillota'
$ws.Cells.Item(5, 16).Value = 1000

# Row 6
$ws.Cells.Item(6, 4).Value = 44167
$ws.Cells.Item(6, 9).Value = 'Primera'
$ws.Cells.Item(6, 10).Value = 1430
$ws.Cells.Item(6, 11).Value = 1200
$ws.Cells.Item(6, 12).Value = 1300
$ws.Cells.Item(6, 13).Value = 1248
$ws.Cells.Item(6, 15).Value = 'Provincia de Quillota'
$ws.Cells.Item(6, 16).Value = 1248

# Row 7
$ws.Cells.Item(7, 4).Value = 44167
$ws.Cells.Item(7, 9).Value = 'Segunda'
$ws.Cells.Item(7, 10).Value = 350
$ws.Cells.Item(7, 11).Value = 1000
$ws.Cells.Item(7, 12).Value = 1000
$ws.Cells.Item(7, 13).Value = 1000
$ws.Cells.Item(7, 15).Value = 'Provincia de Quillota'
$ws.Cells.Item(7, 16).Value = 1000

# Row 8
$ws.Cells.Item(8, 4).Value = 44176
$ws.Cells.Item(8, 9).Value = 'Primera'
$ws.Cells.Item(8, 10).Value = 2500
$ws.Cells.Item(8, 11).Value = 1200
$ws.Cells.Item(8, 12).Value = 1300
$ws.Cells.Item(8, 13).Value = 1256
$ws.Cells.Item(8, 15).Value = 'Provincia de Quillota'
$ws.Cells.Item(8, 16).Value = 1256

# Row 9
$ws.Cells.Item(9, 4).Value = 44176
$ws.Cells.Item(9, 9).Value = 'Segunda'
$ws.Cells.Item(9, 10).Value = 1500
$ws.Cells.Item(9, 11).Value = 1000
$ws.Cells.Item(9, 12).Value = 1000
$ws.Cells.Item(9, 13).Value = 1000
$ws.Cells.Item(9, 15).Value = 'Provincia de Quillota'
$ws.Cells.Item(9, 16).Value = 1000

# Row 10
$ws.Cells.Item(10, 4).Value = 44161
$ws.Cells.Item(10, 9).Value = 'Primera'
$ws.Cells.Item(10, 10).Value = 1600
$ws.Cells.Item(10, 11).Value = 1300
$ws.Cells.Item(10, 12).Value = 1300
$ws.Cells.Item(10, 13).Value = 1300
$ws.Cells.Item(10, 15).Value = 'Provincia de Quillota'
$ws.Cells.Item(10, 16).Value = 1300

# Row 11
$ws.Cells.Item(11, 4).Value = 44161
$ws.Cells.Item(11, 9).Value = 'Segunda'
$ws.Cells.Item(11, 10).Value = 1850
$ws.Cells.Item(11, 11).Value = 1000
$ws.Cells.Item(11, 12).Value = 1000
$ws.Cells.Item(11, 13).Value = 1000
$ws.Cells.Item(11, 15).Value = 'Provincia de Quillota'
$ws.Cells.Item(11, 16).Value = 1000

# Row 12
$ws.Cells.Item(12, 4).Value = 44162
$ws.Cells.Item(12, 9).Value = 'Primera'
$ws.Cells.Item(12, 10).Value = 1200
$ws.Cells.Item(12, 11).Value = 1300
$ws.Cells.Item(12, 12).Value = 1300
$ws.Cells.Item(12, 13).Value = 1300
$ws.Cells.Item(12, 15).Value = 'Provincia de Quillota'
$ws.Cells.Item(12, 16).Value = 1300

# Row 13
$ws.Cells.Item(13, 4).Value = 44162
$ws.Cells.Item(13, 9).Value = 'Segunda'
$ws.Cells.Item(13, 10).Value = 800
$ws.Cells.Item(13, 11).Value = 1000
$ws.Cells.Item(13, 12).Value = 1000
$ws.Cells.Item(13, 13).Value = 1000
$ws.Cells.Item(13, 15).Value = 'Provincia de Quillota'
$ws.Cells.Item(13, 16).Value = 1000

# Row 14
$ws.Cells.Item(14, 4).Value = 44159
$ws.Cells.Item(14, 9).Value = 'Primera'
$ws.Cells.Item(14, 10).Value = 1100
$ws.Cells.Item(14, 11).Value = 1300
$ws.Cells.Item(14, 12).Value = 1300
$ws.Cells.Item(14, 13).Value = 1300
$ws.Cells.Item(14, 15).Value = 'Provincia de Quillota'
$ws.Cells.Item(14, 16).Value = 1300

# Row 15
$ws.Cells.Item(15, 4).Value = 44159
$ws.Cells.Item(15, 9).Value = 'Segunda'
$ws.Cells.Item(15, 10).Value = 800
$ws.Cells.Item(15, 11).Value = 1000
$ws.Cells.Item(15, 12).Value = 1000
$ws.Cells.Item(15, 13).Value = 1000
$ws.Cells.Item(15, 15).Value = 'Provincia de Quillota'
$ws.Cells.Item(15, 16).Value = 1000

# Row 16
$ws.Cells.Item(16, 4).Value = 44169
$ws.Cells.Item(16, 9).Value = 'Primera'
$ws.Cells.Item(16, 10).Value = 950
$ws.Cells.Item(16, 11).Value = 1300
$ws.Cells.Item(16, 12).Value = 1300
$ws.Cells.Item(16, 13).Value = 1300
$ws.Cells.Item(16, 15).Value = 'Provincia de Quillota'
$ws.Cells.Item(16, 16).Value = 1300

# Row 17
$ws.Cells.Item(17, 4).Value = 44169
$ws.Cells.Item(17, 9).Value = 'Segunda'
$ws.Cells.Item(17, 10).Value = 800
$ws.Cells.Item(17, 11).Value = 1000
$ws.Cells.Item(17, 12).Value = 1000
$ws.Cells.Item(17, 13).Value = 1000
$ws.Cells.Item(17, 15).Value = 'Provincia de Quillota'
$ws.Cells.Item(17, 16).Value = 1000

# Row 18
$ws.Cells.Item(18, 4).Value = 44160
$ws.Cells.Item(18, 9).Value = 'Primera'
$ws.Cells.Item(18, 10).Value = 750
$ws.Cells.Item(18, 11).Value = 1300
$ws.Cells.Item(18, 12).Value = 1300
$ws.Cells.Item(18, 13).Value = 1300
$ws.Cells.Item(18, 15).Value = 'Provincia de Quillota'
$ws.Cells.Item(18, 16).Value = 1300

# Row 19
$ws.Cells.Item(19, 4).Value = 44160
$ws.Cells.Item(19, 9).Value = 'Segunda'
$ws.Cells.Item(19, 10).Value = 850
$ws.Cells.Item(19, 11).Value = 1000
$ws.Cells.Item(19, 12).Value = 1000
$ws.Cells.Item(19, 13).Value = 1000
$ws.Cells.Item(19, 15).Value = 'Provincia de Quillota'
$ws.Cells.Item(19, 16).Value = 1000

# Row 20
$ws.Cells.Item(20, 4).Value = 44172
$ws.Cells.Item(20, 9).Value = 'Primera'
$ws.Cells.Item(20, 10).Value = 600
$ws.Cells.Item(20, 11).Value = 1300
$ws.Cells.Item(20, 12).Value = 1300
$ws.Cells.Item(20, 13).Value = 1300
$ws.Cells.Item(20, 15).Value = 'Provincia de Quillota'
$ws.Cells.Item(20, 16).Value = 1300

# Row 21
$ws.Cells.Item(21, 4).Value = 44172
$ws.Cells.Item(21, 9).Value = 'Segunda'
$ws.Cells.Item(21, 10).Value = 550
$ws.Cells.Item(21, 11).Value = 1000
$ws.Cells.Item(21, 12).Value = 1000
$ws.Cells.Item(21, 13).Value = 1000
$ws.Cells.Item(21, 15).Value = 'Provincia de Quillota'
$ws.Cells.Item(21, 16).Value = 1000

# Row 22
$ws.Cells.Item(22, 4).Value = 44168
$ws.Cells.Item(22, 9).Value = 'Primera'
$ws.Cells.Item(22, 10).Value = 1200
$ws.Cells.Item(22, 11).Value = 1300
$ws.Cells.Item(22, 12).Value = 1300
$ws.Cells.Item(22, 13).Value = 1300
$ws.Cells.Item(22, 15).Value = 'Provincia de Quillota'
$ws.Cells.Item(22, 16).Value = 1300

# Row 23
$ws.Cells.Item(23, 4).Value = 44168
$ws.Cells.Item(23, 9).Value = 'Segunda'
$ws.Cells.Item(23, 10).Value = 850
$ws.Cells.Item(23, 11).Value = 1000
$ws.Cells.Item(23, 12).Value = 1000
$ws.Cells.Item(23, 13).Value = 1000
$ws.Cells.Item(23, 15).Value = 'Provincia de Quillota'
$ws.Cells.Item(23, 16).Value = 1000

# Row 24
$ws.Cells.Item(24, 4).Value = 44175
$ws.Cells.Item(24, 9).Value = 'Primera'
$ws.Cells.Item(24, 10).Value = 1500
$ws.Cells.Item(24, 11).Value = 1300
$ws.Cells.Item(24, 12).Value = 1300
$ws.Cells.Item(24, 13).Value = 1300
$ws.Cells.Item(24, 15).Value = 'Provincia de Quillota'
$ws.Cells.Item(24, 16).Value = 1300

# Row 25
$ws.Cells.Item(25, 4).Value = 44175
$ws.Cells.Item(25, 9).Value = 'Segunda'
$ws.Cells.Item(25, 10).Value = 1450
$ws.Cells.Item(25, 11).Value = 1000
$ws.Cells.Item(25, 12).Value = 1000
$ws.Cells.Item(25, 13).Value = 1000
$ws.Cells.Item(25, 15).Value = 'Provincia de Quillota'
$ws.Cells.Item(25, 16).Value = 1000

# Row 26
$ws.Cells.Item(26, 4).Value = 44181
$ws.Cells.Item(26, 9).Value = 'Primera'
$ws.Cells.Item(26, 10).Value = 1000
$ws.Cells.Item(26, 11).Value = 1300
$ws.Cells.Item(26, 12).Value = 1300
$ws.Cells.Item(26, 13).Value = 1300
$ws.Cells.Item(26, 15).Value = 'Provincia de Quillota'
$ws.Cells.Item(26, 16).Value = 1300

# Row 27
$ws.Cells.Item(27, 4).Value = 44181
$ws.Cells.Item(27, 9).Value = 'Segunda'
$ws.Cells.Item(27, 10).Value = 900
$ws.Cells.Item(27, 11).Value = 900
$ws.Cells.Item(27, 12).Value = 900
$ws.Cells.Item(27, 13).Value = 900
$ws.Cells.Item(27, 15).Value = 'Provincia de Quillota'
$ws.Cells.Item(27, 16).Value = 900

# Row 28
$ws.Cells.Item(28, 4).Value = 44179
$ws.Cells.Item(28, 9).Value = 'Primera'
$ws.Cells.Item(28, 10).Value = 980
$ws.Cells.Item(28, 11).Value = 1200
$ws.Cells.Item(28, 12).Value = 1200
$ws.Cells.Item(28, 13).Value = 1200
$ws.Cells.Item(28, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(28, 16).Value = 1200
